$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Kris Dunn", "PG,SG", "LA Clippers"),
    @("Andrew Nembhard", "PG,SG", "Indiana Pacers"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Kyle Kuzma", "SF,PF", "Milwaukee Bucks"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Rows("18:18").Delete()
